$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 35, shifting rows 35:171 down to 36:172.
$ws.Rows.Item(35).Insert()

# The row that used to be row 35 is now row 36; copy its "category" values
# into the brand-new row 35, then overwrite the cells that actually changed.
$src = 36
$dst = 35

$ws.Cells.Item($dst, 1).Value2 = $ws.Cells.Item($src, 1).Value2
$ws.Cells.Item($dst, 2).Value2 = $ws.Cells.Item($src, 2).Value2
$ws.Cells.Item($dst, 3).Value2 = $ws.Cells.Item($src, 3).Value2
$ws.Cells.Item($dst, 4).Value2 = 44453
$ws.Cells.Item($dst, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
$ws.Cells.Item($dst, 5).Value2 = $ws.Cells.Item($src, 5).Value2
$ws.Cells.Item($dst, 6).Value2 = $ws.Cells.Item($src, 6).Value2
$ws.Cells.Item($dst, 7).Value2 = $ws.Cells.Item($src, 7).Value2
$ws.Cells.Item($dst, 8).Value2 = $ws.Cells.Item($src, 8).Value2
$ws.Cells.Item($dst, 9).Value2 = $ws.Cells.Item($src, 9).Value2
$ws.Cells.Item($dst, 10).Value2 = 140
$ws.Cells.Item($dst, 11).Value2 = 9000
$ws.Cells.Item($dst, 12).Value2 = 10000
$ws.Cells.Item($dst, 13).Value2 = 9571
$ws.Cells.Item($dst, 14).Value2 = $ws.Cells.Item($src, 14).Value2
$ws.Cells.Item($dst, 15).Value2 = $ws.Cells.Item($src, 15).Value2
$ws.Cells.Item($dst, 16).Value2 = 1595
$ws.Cells.Item($dst, 17).Value2 = $ws.Cells.Item($src, 17).Value2
$ws.Cells.Item($dst, 18).Value2 = $ws.Cells.Item($src, 18).Value2

# Append a brand-new last row (172) that duplicates the original last
# row's (171, now 172) data -- date 2021-07-23 / 50 units / etc.
$ws.Cells.Item(172, 1).Value2 = 10
$ws.Cells.Item(172, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(172, 3).Value2 = "La Araucanía"
$ws.Cells.Item(172, 4).Value2 = 44400
$ws.Cells.Item(172, 4).NumberFormat = $ws.Cells.Item(171, 4).NumberFormat
$ws.Cells.Item(172, 5).Value2 = 9
$ws.Cells.Item(172, 6).Value2 = 100112017
$ws.Cells.Item(172, 7).Value2 = "Apio"
$ws.Cells.Item(172, 8).Value2 = "Americana (o)"
$ws.Cells.Item(172, 9).Value2 = "Primera"
$ws.Cells.Item(172, 10).Value2 = 50
$ws.Cells.Item(172, 11).Value2 = 9000
$ws.Cells.Item(172, 12).Value2 = 9000
$ws.Cells.Item(172, 13).Value2 = 9000
$ws.Cells.Item(172, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(172, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(172, 16).Value2 = 1500
$ws.Cells.Item(172, 17).Value2 = 6
$ws.Cells.Item(172, 18).Value2 = "Hortaliza"
